$wb = $excel.ActiveWorkbook

$oldGuid = "a1f78878-6f34-4aba-8c6b-09ecfc950b78"
$newGuid = "634468d4-c794-4a09-b122-e59c14183801"

$newZhHash = "d0dffd5941a32359304aa6704160ef35057c3a98"

function Set-HyperlinkDisplay($ws, $address, $text) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $address) {
            $hl.TextToDisplay = $text
        }
    }
}

function Remove-HyperlinkAt($ws, $address) {
    $toDelete = @()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $address) {
            $toDelete += $hl
        }
    }
    foreach ($hl in $toDelete) {
        $hl.Delete()
    }
}

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
Set-HyperlinkDisplay $wsOverview '$B$2' "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-27 06:56:33"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newGuid.md"
Set-HyperlinkDisplay $wsZh '$A$2' "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newZhHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-27 06:56:29"
Remove-HyperlinkAt $wsZh '$I$2'
$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newGuid.md"
Set-HyperlinkDisplay $wsDe '$A$2' "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newZhHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-27 06:56:33"
Remove-HyperlinkAt $wsDe '$I$2'
$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"
